# Auto-generated Excel COM-interop script applying the recorded diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33
$ws.Range("A33").Value = 111959819
$ws.Range("B33").Value = 89834
$ws.Range("E33").Value = 658
$ws.Range("F33").Value = "Rosenticka"
$ws.Range("G33").Value = "Rhodofomes roseus"
$ws.Range("H33").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("K33").Value = ""
$ws.Range("L33").Value = ""
$ws.Range("M33").Value = ""
$ws.Range("N33").Value = ""
$ws.Range("Q33").Value = 561895
$ws.Range("R33").Value = 7307265

# Row 34
$ws.Range("A34").Value = 111959829
$ws.Range("B34").Value = 77650
$ws.Range("E34").Value = 6425
$ws.Range("F34").Value = "Garnlav"
$ws.Range("G34").Value = "Alectoria sarmentosa"
$ws.Range("H34").Value = "(Ach.) Ach."
$ws.Range("Q34").Value = 561894
$ws.Range("R34").Value = 7307220

# Row 35
$ws.Range("A35").Value = 111985428
$ws.Range("B35").Value = 56430
$ws.Range("E35").Value = 100109
$ws.Range("F35").Value = "Tretåig hackspett"
$ws.Range("G35").Value = "Picoides tridactylus"
$ws.Range("H35").Value = "(Linnaeus, 1758)"
$ws.Range("K35").Value = ""
$ws.Range("L35").Value = ""
$ws.Range("M35").Value = "äldre spår"
$ws.Range("N35").Value = ""
$ws.Range("Q35").Value = 561916
$ws.Range("R35").Value = 7307339

# Row 45
$ws.Range("A45").Value = 111985425
$ws.Range("B45").Value = 56430
$ws.Range("E45").Value = 100109
$ws.Range("F45").Value = "Tretåig hackspett"
$ws.Range("G45").Value = "Picoides tridactylus"
$ws.Range("H45").Value = "(Linnaeus, 1758)"
$ws.Range("K45").Value = ""
$ws.Range("L45").Value = ""
$ws.Range("M45").Value = "färska spår"
$ws.Range("N45").Value = ""
$ws.Range("Q45").Value = 561986
$ws.Range("R45").Value = 7307363

# Row 46
$ws.Range("A46").Value = 111959827
$ws.Range("B46").Value = 89571
$ws.Range("E46").Value = 5432
$ws.Range("F46").Value = "Granticka"
$ws.Range("G46").Value = "Porodaedalea chrysoloma"
$ws.Range("H46").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q46").Value = 561870
$ws.Range("R46").Value = 7307209

# Row 47
$ws.Range("A47").Value = 111959821
$ws.Range("B47").Value = 77731
$ws.Range("E47").Value = 864
$ws.Range("F47").Value = "Knottrig blåslav"
$ws.Range("G47").Value = "Hypogymnia bitteri"
$ws.Range("H47").Value = "(Lynge) Ahti"
$ws.Range("K47").Value = ""
$ws.Range("L47").Value = ""
$ws.Range("M47").Value = ""
$ws.Range("N47").Value = ""
$ws.Range("Q47").Value = 561962
$ws.Range("R47").Value = 7307351

# Row 86
$ws.Range("A86").Value = 112045886
$ws.Range("B86").Value = 56430
$ws.Range("D86").Value = "NT"
$ws.Range("E86").Value = 100109
$ws.Range("F86").Value = "Tretåig hackspett"
$ws.Range("G86").Value = "Picoides tridactylus"
$ws.Range("H86").Value = "(Linnaeus, 1758)"
$ws.Range("M86").Value = "färska spår"
$ws.Range("Q86").Value = 561826
$ws.Range("R86").Value = 7307246
$ws.Range("Z86").Value = "11:50"
$ws.Range("AB86").Value = "11:50"

# Row 87
$ws.Range("A87").Value = 112045880
$ws.Range("B87").Value = 89993
$ws.Range("D87").Value = "VU"
$ws.Range("E87").Value = 1209
$ws.Range("F87").Value = "Rynkskinn"
$ws.Range("G87").Value = "Phlebia centrifuga"
$ws.Range("H87").Value = "P.Karst."
$ws.Range("M87").Value = ""
$ws.Range("Q87").Value = 561969
$ws.Range("R87").Value = 7307275
$ws.Range("Z87").Value = "12:57"
$ws.Range("AB87").Value = "12:57"

# Row 108
$ws.Range("A108").Value = 112045898
$ws.Range("B108").Value = 56430
$ws.Range("E108").Value = 100109
$ws.Range("F108").Value = "Tretåig hackspett"
$ws.Range("G108").Value = "Picoides tridactylus"
$ws.Range("H108").Value = "(Linnaeus, 1758)"
$ws.Range("M108").Value = "färska spår"
$ws.Range("Q108").Value = 561992
$ws.Range("R108").Value = 7307388
$ws.Range("Z108").Value = "10:12"
$ws.Range("AB108").Value = "10:12"

# Row 109
$ws.Range("A109").Value = 112045897
$ws.Range("B109").Value = 73834
$ws.Range("E109").Value = 6440
$ws.Range("F109").Value = "Vitgrynig nållav"
$ws.Range("G109").Value = "Chaenotheca subroscida"
$ws.Range("H109").Value = "(Eitner) Zahlbr."
$ws.Range("M109").Value = ""
$ws.Range("Q109").Value = 561955
$ws.Range("R109").Value = 7307352
$ws.Range("Z109").Value = "10:34"
$ws.Range("AB109").Value = "10:34"

# Row 110
$ws.Range("A110").Value = 112045878
$ws.Range("B110").Value = 89553
$ws.Range("E110").Value = 1202
$ws.Range("F110").Value = "Ullticka"
$ws.Range("G110").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H110").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q110").Value = 561966
$ws.Range("R110").Value = 7307274
$ws.Range("Z110").Value = "12:57"
$ws.Range("AB110").Value = "12:57"

# Row 111
$ws.Range("A111").Value = 112045876
$ws.Range("B111").Value = 77650
$ws.Range("E111").Value = 6425
$ws.Range("F111").Value = "Garnlav"
$ws.Range("G111").Value = "Alectoria sarmentosa"
$ws.Range("H111").Value = "(Ach.) Ach."
$ws.Range("Q111").Value = 561966
$ws.Range("R111").Value = 7307274
$ws.Range("Z111").Value = "12:58"
$ws.Range("AB111").Value = "12:58"

